$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 95

$ws.Cells.Item($newRow, 1).Value = 46044
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = 221
$ws.Cells.Item($newRow, 3).Value = 227
$ws.Cells.Item($newRow, 4).Value = 215
